# appendix/metrics_all_sentiment-news-econ.xlsx
# Commit: "added balanced accuracy based on reviewer feedback"
#
# The sheet currently has two 6-row metric blocks after the "mean" blocks:
#   rows 14-19: f1_macro_std
#   rows 20-25: f1_micro_std
#
# The new layout inserts a new "accuracy_balanced_mean" block right after the
# f1_micro_mean block (so it becomes rows 14-19), which pushes the existing
# f1_macro_std / f1_micro_std blocks down by 6 rows (to rows 20-25 / 26-31),
# and appends a new "accuracy_balanced_std" block at the bottom (rows 32-37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert 6 blank rows at row 14 for accuracy_balanced_mean.
#     This shifts the existing f1_macro_std (old 14-19) / f1_micro_std (old
#     20-25) blocks down to rows 20-25 / 26-31 -- their label, value and
#     style content travels with them, untouched.
$ws.Range("A14:H19").EntireRow.Insert()

# --- 2. Insert 6 more blank rows at row 32 (right after the shifted
#     f1_micro_std block) to hold the new accuracy_balanced_std block.
$ws.Range("A32:H37").EntireRow.Insert()

# --- 3. Newly inserted rows don't inherit the bold/bordered/centered style
#     used for the column-A metric-name cells, so copy that formatting from
#     an existing labeled cell (row 13, still "f1_micro_mean" / style intact).
$ws.Range("A13").Copy()
$ws.Range("A14:A19").PasteSpecial(-4122)
$ws.Range("A32:A37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 4. Row/column headers for the two new blocks.
$ws.Range("A14:A19").Value = "accuracy_balanced_mean"
$ws.Range("A32:A37").Value = "accuracy_balanced_std"

# --- 5. accuracy_balanced_mean data (rows 14-19).
$ws.Cells.Item(14, 2).Value = "'0"
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0.712
$ws.Cells.Item(15, 2).Value = "'100"
$ws.Cells.Item(15, 3).Value = 0.534
$ws.Cells.Item(15, 4).Value = 0.53
$ws.Cells.Item(15, 5).Value = 0.564
$ws.Cells.Item(15, 6).Value = 0.551
$ws.Cells.Item(15, 7).Value = 0.574
$ws.Cells.Item(15, 8).Value = 0.677
$ws.Cells.Item(16, 2).Value = "'500"
$ws.Cells.Item(16, 3).Value = 0.574
$ws.Cells.Item(16, 4).Value = 0.592
$ws.Cells.Item(16, 5).Value = 0.616
$ws.Cells.Item(16, 6).Value = 0.558
$ws.Cells.Item(16, 7).Value = 0.691
$ws.Cells.Item(16, 8).Value = 0.703
$ws.Cells.Item(17, 2).Value = "'1000"
$ws.Cells.Item(17, 3).Value = 0.602
$ws.Cells.Item(17, 4).Value = 0.613
$ws.Cells.Item(17, 5).Value = 0.603
$ws.Cells.Item(17, 6).Value = 0.581
$ws.Cells.Item(17, 7).Value = 0.712
$ws.Cells.Item(17, 8).Value = 0.742
$ws.Cells.Item(18, 2).Value = "'2500"
$ws.Cells.Item(18, 3).Value = 0.649
$ws.Cells.Item(18, 4).Value = 0.637
$ws.Cells.Item(18, 5).Value = 0.612
$ws.Cells.Item(18, 6).Value = 0.608
$ws.Cells.Item(18, 7).Value = 0.695
$ws.Cells.Item(18, 8).Value = 0.747
$ws.Cells.Item(19, 2).Value = "3000 (all)"
$ws.Cells.Item(19, 3).Value = 0.669
$ws.Cells.Item(19, 4).Value = 0.65
$ws.Cells.Item(19, 5).Value = 0.619
$ws.Cells.Item(19, 6).Value = 0.595
$ws.Cells.Item(19, 7).Value = 0.702
$ws.Cells.Item(19, 8).Value = 0.739

# --- 6. accuracy_balanced_std data (rows 32-37).
$ws.Cells.Item(32, 2).Value = "'0"
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(33, 2).Value = "'100"
$ws.Cells.Item(33, 3).Value = 0.034
$ws.Cells.Item(33, 4).Value = 0.032
$ws.Cells.Item(33, 5).Value = 0.06
$ws.Cells.Item(33, 6).Value = 0.04
$ws.Cells.Item(33, 7).Value = 0.079
$ws.Cells.Item(33, 8).Value = 0.027
$ws.Cells.Item(34, 2).Value = "'500"
$ws.Cells.Item(34, 3).Value = 0.028
$ws.Cells.Item(34, 4).Value = 0.034
$ws.Cells.Item(34, 5).Value = 0.03
$ws.Cells.Item(34, 6).Value = 0.046
$ws.Cells.Item(34, 7).Value = 0.018
$ws.Cells.Item(34, 8).Value = 0.016
$ws.Cells.Item(35, 2).Value = "'1000"
$ws.Cells.Item(35, 3).Value = 0.014
$ws.Cells.Item(35, 4).Value = 0.023
$ws.Cells.Item(35, 5).Value = 0.015
$ws.Cells.Item(35, 6).Value = 0.029
$ws.Cells.Item(35, 7).Value = 0.009
$ws.Cells.Item(35, 8).Value = 0.01
$ws.Cells.Item(36, 2).Value = "'2500"
$ws.Cells.Item(36, 3).Value = 0.004
$ws.Cells.Item(36, 4).Value = 0.002
$ws.Cells.Item(36, 5).Value = 0.011
$ws.Cells.Item(36, 6).Value = 0.005
$ws.Cells.Item(36, 7).Value = 0.019
$ws.Cells.Item(36, 8).Value = 0.009
$ws.Cells.Item(37, 2).Value = "3000 (all)"
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(37, 4).Value = 0.004
$ws.Cells.Item(37, 5).Value = 0
$ws.Cells.Item(37, 6).Value = 0.011
$ws.Cells.Item(37, 7).Value = 0.015
$ws.Cells.Item(37, 8).Value = 0.005
